$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force Excel to store the value as literal text, matching the
    # original inline-string cell type, even for digit/period-only
    # strings that Excel would otherwise auto-convert to a number.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Updates where Excel will not misinterpret the string as a number ---
$ws.Range("D2").Value = "27.203.64"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.905.32"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  +10.34%  "
$ws.Range("D13").Value = "1.914.76"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "27.242.64"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "2.158.85"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("E38").Value = "  +3.45%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -0.04%  "

# --- Updates where the text looks numeric and must be forced to text ---
Set-TextValue "D5" "307.85"
Set-TextValue "D6" "1.001"
Set-TextValue "D7" "0.5200"
Set-TextValue "D8" "0.3767"
Set-TextValue "D9" "0.07273"
Set-TextValue "D12" "0.08453"
Set-TextValue "D14" "96.86"
Set-TextValue "D15" "5.299"
Set-TextValue "D16" "1.003"
Set-TextValue "D17" "0.000008666"
Set-TextValue "D18" "14.53"
Set-TextValue "D21" "5.098"
Set-TextValue "D23" "10.65"
Set-TextValue "D24" "6.444"
Set-TextValue "D25" "2.343"
Set-TextValue "D26" "146.89"
Set-TextValue "D27" "1.757"
Set-TextValue "D28" "18.25"
Set-TextValue "D29" "115.19"
Set-TextValue "D30" "4.823"
Set-TextValue "D31" "4.908"
Set-TextValue "D32" "0.09281"
Set-TextValue "D33" "0.05079"
Set-TextValue "D34" "0.7955"
Set-TextValue "D36" "3.426"
Set-TextValue "D37" "2.948"
Set-TextValue "D38" "0.5820"
Set-TextValue "D40" "0.02008"
Set-TextValue "D42" "9.070"
Set-TextValue "D44" "116.77"
Set-TextValue "D46" "0.4895"
Set-TextValue "D47" "1.001"
Set-TextValue "D48" "10.17"
Set-TextValue "D49" "1.639"
Set-TextValue "D50" "37.73"
Set-TextValue "D51" "64.10"
